$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number and must be forced to stay
# text (matching the workbook's convention of storing prices as strings).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

function Set-PlainValue($cellRef, $text) {
    $ws.Range($cellRef).Value = $text
}

# Row 2 - Bitcoin
Set-PlainValue "D2" "26.732.48"
Set-PlainValue "E2" "  +0.24%  "

# Row 3 - Ethereum
Set-PlainValue "D3" "1.602.49"

# Row 5 - BNB
Set-TextValue "D5" "211.75"
Set-PlainValue "E5" "  -0.04%  "

# Row 6 - XRP
Set-TextValue "D6" "0.512"
Set-PlainValue "E6" "  -0.58%  "

# Row 7 - USDC
Set-PlainValue "E7" "  +0.19%  "

# Row 8 - Dogecoin
Set-PlainValue "E8" "  +0.21%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.247"
Set-PlainValue "E9" "  +0.22%  "

# Row 10 - Solana
Set-TextValue "D10" "19.73"
Set-PlainValue "E10" "  +0.83%  "

# Row 11 - TRON
Set-PlainValue "E11" "  +1.00%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-PlainValue "D12" "1.827.51"
Set-PlainValue "E12" "  +0.22%  "

# Row 13 - WrappedEther
Set-PlainValue "D13" "1.614.73"
Set-PlainValue "E13" "  +2.00%  "

# Row 14 - Polkadot
Set-PlainValue "E14" "  +0.42%  "

# Row 15 - Polygon
Set-PlainValue "E15" "  -0.49%  "

# Row 16 - Litecoin
Set-TextValue "D16" "65.06"
Set-PlainValue "E16" "  -0.24%  "

# Row 17 - ShibaInu
Set-PlainValue "E17" "  +0.39%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "210.34"
Set-PlainValue "E18" "  +0.57%  "

# Row 19 & 20 - swap Chainlink/Dai with updated values
Set-PlainValue "B19" "Dai"
Set-PlainValue "C19" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D19" "1.01"
Set-PlainValue "E19" "  +0.19%  "

Set-PlainValue "B20" "Chainlink"
Set-PlainValue "C20" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D20" "7.18"
Set-PlainValue "E20" "  +1.95%  "

# Row 21 - Uniswap
Set-TextValue "D21" "4.28"
Set-PlainValue "E21" "  -0.19%  "

# Row 22 - Toncoin
Set-TextValue "D22" "2.27"
Set-PlainValue "E22" "  -2.01%  "

# Row 24 - Monero
Set-TextValue "D24" "143.62"
Set-PlainValue "E24" "  -1.18%  "

# Row 25 - BinanceUSD
Set-PlainValue "E25" "  +0.20%  "

# Row 26 - Cosmos
Set-TextValue "D26" "7.09"
Set-PlainValue "E26" "  -0.24%  "

# Row 27 - Stellar
Set-PlainValue "E27" "  -0.85%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "15.39"
Set-PlainValue "E28" "  +0.65%  "

# Row 29 - Hedera
Set-PlainValue "E29" "  -0.69%  "

# Row 30 - PancakeSwap
Set-PlainValue "E30" "  -0.11%  "

# Row 31 - Filecoin
Set-PlainValue "E31" "  +1.29%  "

# Row 32 - InternetComputer(DFINITY)
Set-PlainValue "E32" "  +1.05%  "

# Row 33 - Maker
Set-PlainValue "D33" "1.294.07"
Set-PlainValue "E33" "  +1.13%  "

# Row 34 - HuobiToken
Set-PlainValue "E34" "  +0.77%  "

# Row 35 - LidoDAOToken
Set-PlainValue "E35" "  +0.74%  "

# Row 36 - ImmutableX
Set-PlainValue "E36" "  -3.39%  "

# Row 37 - WEMIXToken
Set-PlainValue "E37" "  +10.70%  "

# Row 38 - VeChain
Set-TextValue "D38" "0.0169"
Set-PlainValue "E38" "  -0.34%  "

# Row 40 - FraxShare
Set-PlainValue "E40" "  -1.89%  "

# Row 41 - MXToken
Set-PlainValue "E41" "  -0.32%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "0.785"
Set-PlainValue "E42" "  +0.18%  "

# Row 43 - Aave
Set-TextValue "D43" "63.04"
Set-PlainValue "E43" "  -1.42%  "

# Row 44 - RocketPoolETH
Set-PlainValue "D44" "1.738.88"
Set-PlainValue "E44" "  +0.16%  "

# Row 45 - Quant
Set-TextValue "D45" "90.61"
Set-PlainValue "E45" "  -0.79%  "

# Row 46 - RenderToken
Set-PlainValue "E46" "  -2.38%  "

# Row 47 - Algorand
Set-PlainValue "E47" "  -0.63%  "

# Row 48 - Cronos
Set-PlainValue "E48" "  +1.66%  "

# Row 49 - USDD
Set-PlainValue "E49" "  +0.21%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "7.43"
Set-PlainValue "E50" "  +0.26%  "
